# Add a new "2021" column (R) to the tourism-GDP-share table, mirroring the
# formatting already used by the preceding "2020" column (Q), and move the
# active selection to O9 (as recorded by the author's session after the
# edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (years header): R4 = 2021, formatted like Q4 -------------------
$ws.Range("Q4").Copy() | Out-Null
$ws.Range("R4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("R4").Value = 2021

# --- Row 5 (data row): R5 = 3.6, formatted like Q5 -------------------------
$ws.Range("Q5").Copy() | Out-Null
$ws.Range("R5").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("R5").Value = 3.6

$excel.CutCopyMode = $false

# --- Move the selection to match the saved view ----------------------------
$ws.Range("O9").Select() | Out-Null
